$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.24
$ws.Range("X2").Value = 970
$ws.Range("Y2").Value = 970
$ws.Range("AB2").Value = 970
$ws.Range("AC2").Value = 970
$ws.Range("AD2").Value = 970
$ws.Range("AG2").Value = 970
$ws.Range("AH2").Value = 970
$ws.Range("G4").Value = 1.82
$ws.Range("W4").Value = 2.2
$ws.Range("AL4").Value = 50
$ws.Range("AN4").Value = 14.5
$ws.Range("G5").Value = 1.79
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 6.6
$ws.Range("J5").Value = 3.65
$ws.Range("K5").Value = 4.1
$ws.Range("N5").Value = 3.3
$ws.Range("O5").Value = 1.36
$ws.Range("R5").Value = 1.3
$ws.Range("V5").Value = 1.18
$ws.Range("W5").Value = 2.26
$ws.Range("AB5").Value = 8.4
$ws.Range("AG5").Value = 1000
$ws.Range("H6").Value = 14
$ws.Range("J6").Value = 7.2
$ws.Range("N6").Value = 7
$ws.Range("O6").Value = 1.13
$ws.Range("P6").Value = 3.05
$ws.Range("Q6").Value = 1.35
$ws.Range("R6").Value = 1.82
$ws.Range("S6").Value = 1.97
$ws.Range("T6").Value = 1.94
$ws.Range("U6").Value = 1.86
$ws.Range("W6").Value = 4.8
$ws.Range("X6").Value = 46
$ws.Range("Y6").Value = 970
$ws.Range("Z6").Value = 180
$ws.Range("AB6").Value = 15
$ws.Range("AC6").Value = 22
$ws.Range("AD6").Value = 970
$ws.Range("AE6").Value = 250
$ws.Range("AF6").Value = 11.5
$ws.Range("AG6").Value = 14
$ws.Range("AH6").Value = 38
$ws.Range("AI6").Value = 170
$ws.Range("AJ6").Value = 11.5
$ws.Range("AK6").Value = 16
$ws.Range("AL6").Value = 42
$ws.Range("AM6").Value = 170
$ws.Range("AN6").Value = 3.7
$ws.Range("K7").Value = 12.5
$ws.Range("X7").Value = 970
$ws.Range("AB7").Value = 970
$ws.Range("AC7").Value = 970
$ws.Range("AG7").Value = 970
$ws.Range("AH7").Value = 970
$ws.Range("T8").Value = 2.4
$ws.Range("U8").Value = 1.6
$ws.Range("H9").Value = 5.7
$ws.Range("G10").Value = 1.89
$ws.Range("N10").Value = 2.62
$ws.Range("S10").Value = 3.45
$ws.Range("T10").Value = 1.83
$ws.Range("U10").Value = 1.68
$ws.Range("Y10").Value = 970
$ws.Range("AD10").Value = 970
$ws.Range("AH10").Value = 970
$ws.Range("AN10").Value = 20
$ws.Range("F11").Value = 2.84
$ws.Range("G11").Value = 3.25
$ws.Range("H11").Value = 2.78
$ws.Range("I11").Value = 3.2
$ws.Range("J11").Value = 2.76
$ws.Range("K11").Value = 3.2
$ws.Range("L11").Value = 1.01
$ws.Range("M11").Value = 1.11
$ws.Range("N11").Value = 2.46
$ws.Range("O11").Value = 1.55
$ws.Range("P11").Value = 1.48
$ws.Range("Q11").Value = 2.62
$ws.Range("R11").Value = 1.17
$ws.Range("S11").Value = 4.9
$ws.Range("T11").Value = 2.1
$ws.Range("U11").Value = 1.75
$ws.Range("V11").Value = 1.46
$ws.Range("W11").Value = 1.44
$ws.Range("X11").Value = 8.4
$ws.Range("Y11").Value = 8.6
$ws.Range("Z11").Value = 19
$ws.Range("AA11").Value = 60
$ws.Range("AB11").Value = 8.800000000000001
$ws.Range("AC11").Value = 7.2
$ws.Range("AD11").Value = 14.5
$ws.Range("AE11").Value = 48
$ws.Range("AF11").Value = 19.5
$ws.Range("AG11").Value = 15
$ws.Range("AH11").Value = 24
$ws.Range("AI11").Value = 80
$ws.Range("AJ11").Value = 60
$ws.Range("AK11").Value = 50
$ws.Range("AL11").Value = 80
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 65
$ws.Range("AO11").Value = 60
$ws.Range("K12").Value = 3.55
